# edit imageid.xlsx for future reference
#
# Renames two of the header labels on Sheet1 and updates the active
# selection to reflect where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: rename "replace_dir" -> "greek_dir" and "shuffle_dir" -> "scramble_dir"
$ws.Range("D1").Value = "greek_dir"
$ws.Range("E1").Value = "scramble_dir"

# Move/save the active cell selection to E2
$ws.Activate()
$ws.Range("E2").Select()
